# Applies commit "subject 17 - 20": adds a "Valid" column (L) to the
# prequestionnaire sheet, fills in previously-blank rows for subjects 12-20
# (rows 14-22), fixes a "Joysick" -> "Joystick" typo in G12, and fills the
# Condition (K) / Valid (L) columns for the remaining subjects (rows 23-32).
#
# NOTE on write order: the xlsx shared-string table is append-only (new
# strings are appended in first-reference order, and strings that drop to
# zero references get garbage collected on save). The target file expects
# brand new strings to land in a specific order, so the statements below are
# deliberately sequenced to reproduce it:
#   1. "Valid"                       (header L1)
#   2. rows 14-19 (left to right)    -> "DeskTop; Tablet; Mobile;", "Biochemistry",
#                                        "HE", "Joystick; "
#   3. G12 typo fix                  -> "Joystick; Wiimote"
#   4. rows 20-22 (left to right)    -> "BME", "N/A"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "Valid" column header (introduces shared string "Valid")
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 12).Value = "Valid"

# ---------------------------------------------------------------------
# 2. Subjects 12-17 (rows 14-19) - previously-blank data rows
# ---------------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "Female"
$ws.Cells.Item(14, 3).Value = 21
$ws.Cells.Item(14, 4).Value = "RBE"
$ws.Cells.Item(14, 5).Value = 4
$ws.Cells.Item(14, 6).Value = "Desktop; Game Console"
$ws.Cells.Item(14, 7).Value = "Keyboard/Mouse; Joystick; "
$ws.Cells.Item(14, 8).Value = 2
$ws.Cells.Item(14, 9).Value = 1
$ws.Cells.Item(14, 10).Value = 4
$ws.Cells.Item(14, 11).Value = "G"
$ws.Cells.Item(14, 12).Value = 1

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Other"
$ws.Cells.Item(15, 3).Value = 22
$ws.Cells.Item(15, 4).Value = "Psychology"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = "Mobile"
$ws.Cells.Item(15, 7).Value = "Other"
$ws.Cells.Item(15, 8).Value = 1
$ws.Cells.Item(15, 9).Value = 1
$ws.Cells.Item(15, 10).Value = 1
$ws.Cells.Item(15, 11).Value = "F"
$ws.Cells.Item(15, 12).Value = 0

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Other"
$ws.Cells.Item(16, 3).Value = 18
$ws.Cells.Item(16, 4).Value = "CS"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = "DeskTop; Tablet; Mobile;"
$ws.Cells.Item(16, 7).Value = "Keyboard/Mouse"
$ws.Cells.Item(16, 8).Value = 2
$ws.Cells.Item(16, 9).Value = 1
$ws.Cells.Item(16, 10).Value = 3
$ws.Cells.Item(16, 11).Value = "G"
$ws.Cells.Item(16, 12).Value = 0

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "Other"
$ws.Cells.Item(17, 3).Value = 21
$ws.Cells.Item(17, 4).Value = "Biochemistry"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = "Desktop"
$ws.Cells.Item(17, 7).Value = "Keyboard/Mouse"
$ws.Cells.Item(17, 8).Value = 1
$ws.Cells.Item(17, 9).Value = 1
$ws.Cells.Item(17, 10).Value = 4
$ws.Cells.Item(17, 11).Value = "F"
$ws.Cells.Item(17, 12).Value = 1

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "Female"
$ws.Cells.Item(18, 3).Value = 22
$ws.Cells.Item(18, 4).Value = "RBE"
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(18, 6).Value = "Desktop"
$ws.Cells.Item(18, 7).Value = "Keyboard/Mouse"
$ws.Cells.Item(18, 8).Value = 4
$ws.Cells.Item(18, 9).Value = 2
$ws.Cells.Item(18, 10).Value = 3
$ws.Cells.Item(18, 11).Value = "G"
$ws.Cells.Item(18, 12).Value = 1

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "Male"
$ws.Cells.Item(19, 3).Value = 21
$ws.Cells.Item(19, 4).Value = "HE"
$ws.Cells.Item(19, 5).Value = 4
$ws.Cells.Item(19, 6).Value = "Game Console"
$ws.Cells.Item(19, 7).Value = "Joystick; "
$ws.Cells.Item(19, 8).Value = 1
$ws.Cells.Item(19, 9).Value = 1
$ws.Cells.Item(19, 10).Value = 3
$ws.Cells.Item(19, 11).Value = "F"
$ws.Cells.Item(19, 12).Value = 1

# ---------------------------------------------------------------------
# 3. Fix "Joysick" typo on row 12 -> "Joystick; Wiimote"
# ---------------------------------------------------------------------
$ws.Cells.Item(12, 7).Value = "Joystick; Wiimote"

# ---------------------------------------------------------------------
# 4. Subjects 18-20 (rows 20-22)
# ---------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "Female"
$ws.Cells.Item(20, 3).Value = 19
$ws.Cells.Item(20, 4).Value = "Biology"
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 6).Value = "Desktop"
$ws.Cells.Item(20, 7).Value = "Keyboard/Mouse"
$ws.Cells.Item(20, 8).Value = 3
$ws.Cells.Item(20, 9).Value = 4
$ws.Cells.Item(20, 10).Value = 2
$ws.Cells.Item(20, 11).Value = "G"
$ws.Cells.Item(20, 12).Value = 1

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "Male"
$ws.Cells.Item(21, 3).Value = 19
$ws.Cells.Item(21, 4).Value = "BME"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = "DeskTop; Tablet; Mobile;"
$ws.Cells.Item(21, 7).Value = "Keyboard/Mouse"
$ws.Cells.Item(21, 8).Value = 3
$ws.Cells.Item(21, 9).Value = 3
$ws.Cells.Item(21, 10).Value = 3
$ws.Cells.Item(21, 11).Value = "F"
$ws.Cells.Item(21, 12).Value = 1

$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "Male"
$ws.Cells.Item(22, 3).Value = "N/A"
$ws.Cells.Item(22, 4).Value = "N/A"
$ws.Cells.Item(22, 5).Value = 5
$ws.Cells.Item(22, 6).Value = "Game Console"
$ws.Cells.Item(22, 7).Value = "Keyboard/Mouse; Joystick; "
$ws.Cells.Item(22, 8).Value = 5
$ws.Cells.Item(22, 9).Value = 3
$ws.Cells.Item(22, 10).Value = 2
$ws.Cells.Item(22, 11).Value = "G"
$ws.Cells.Item(22, 12).Value = 1

# ---------------------------------------------------------------------
# 5. "Valid" column for the already-existing subjects 0-11 (rows 2-13)
#    Row 9 (subject 7) intentionally has no "Valid" value, matching source.
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(13, 12).Value = 1

# ---------------------------------------------------------------------
# 6. Condition (K) / Valid (L) for subjects 21-30 (rows 23-32),
#    alternating "F"/"G" condition, continuing the existing pattern.
# ---------------------------------------------------------------------
$ws.Cells.Item(23, 11).Value = "F"
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(24, 11).Value = "G"
$ws.Cells.Item(24, 12).Value = 1
$ws.Cells.Item(25, 11).Value = "F"
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(26, 11).Value = "G"
$ws.Cells.Item(26, 12).Value = 1
$ws.Cells.Item(27, 11).Value = "F"
$ws.Cells.Item(27, 12).Value = 1
$ws.Cells.Item(28, 11).Value = "G"
$ws.Cells.Item(28, 12).Value = 1
$ws.Cells.Item(29, 11).Value = "F"
$ws.Cells.Item(29, 12).Value = 1
$ws.Cells.Item(30, 11).Value = "G"
$ws.Cells.Item(30, 12).Value = 1
$ws.Cells.Item(31, 11).Value = "F"
$ws.Cells.Item(31, 12).Value = 1
$ws.Cells.Item(32, 11).Value = "G"
$ws.Cells.Item(32, 12).Value = 1

# ---------------------------------------------------------------------
# 7. Update the view: scroll so row 4 is at the top and select K22
#    (matches the sheetView/selection in the target file).
# ---------------------------------------------------------------------
$ws.Range("K22").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
